# Refresh the cryptos list (GitHub Actions style update): new Price (column D)
# and Volume(1h) (column E) values for rows 2-51.
#
# Column D values are stored as text in the workbook (e.g. "51.101.54" uses
# "." as a thousands separator, and trailing zeros like "374.30" must be kept
# literally). Plain decimal-looking strings would otherwise be auto-detected
# as numbers by Excel (dropping the trailing zero / changing type), so each
# D write is entered with a leading apostrophe to force literal text, and the
# cell style is then reset to Normal so the quote-prefix marker doesn't leave
# a lingering formatting change on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'51.037.16"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = "'2.940.60"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'374.30"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('D6').Value = "'101.09"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.36%  '
$ws.Range('D7').Value = "'0.536"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.56%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range('D10').Value = "'36.29"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.73%  '
$ws.Range('D11').Value = "'0.138"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').Value = "'0.0848"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = "'3.408.57"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.47%  '
$ws.Range('D14').Value = "'18.02"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('D15').Value = "'7.56"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').Value = "'2.944.74"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = "'0.993"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.58%  '
$ws.Range('D18').Value = "'10.90"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +45.98%  '
$ws.Range('D19').Value = "'50.992.72"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.04%  '
$ws.Range('D20').Value = "'3.09"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.96%  '
$ws.Range('D21').Value = "'12.40"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.48%  '
$ws.Range('D22').Value = "'0.0₃0955"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.03%  '
$ws.Range('D23').Value = "'264.86"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('D24').Value = "'68.65"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.12%  '
$ws.Range('D25').Value = "'3.13"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.60%  '
$ws.Range('D26').Value = "'8.14"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.71%  '
$ws.Range('D27').Value = "'7.61"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('D28').Value = "'0.999"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = "'25.57"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.63%  '
$ws.Range('D30').Value = "'0.163"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.31%  '
$ws.Range('D31').Value = "'0.109"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.43%  '
$ws.Range('D32').Value = "'10.00"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('D33').Value = "'50.78"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('E34').Value = '  -1.13%  '
$ws.Range('D35').Value = "'33.30"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.88%  '
$ws.Range('D36').Value = "'0.0441"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.64%  '
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('E38').Value = '  +3.59%  '
$ws.Range('E39').Value = '  -1.33%  '
$ws.Range('D40').Value = "'16.31"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.29%  '
$ws.Range('E41').Value = '  -3.22%  '
$ws.Range('E42').Value = '  -4.20%  '
$ws.Range('D43').Value = "'120.47"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.97%  '
$ws.Range('D44').Value = "'21.32"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.15%  '
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('D47').Value = "'0.271"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.97%  '
$ws.Range('D48').Value = "'2.32"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('D49').Value = "'1.988.66"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.17%  '
$ws.Range('D50').Value = "'0.0323"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.43%  '
$ws.Range('E51').Value = '  +1.79%  '
